$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = 80976130
$ws.Range("L62").Value = 'hane'
$ws.Range("M62").Value = 'frispringande/krypande'
$ws.Range("AC62").Value = 'Ny lokal, och blott tredje kända i Södermanland!? Närmast funnen i Tyresta NP. Bör eftersökas på fler lokaler i kommunen! Grävde först fram fragment av en ad hona, men kunde sedan finna en vuxen hane i en perfekt rödmurken granlåga i sent nedbrytningsstadium.'
$ws.Range("A63").Value = 111683856
$ws.Range("B63").Value = 108219
$ws.Range("E63").Value = 219711
$ws.Range("F63").Value = 'Sårläka'
$ws.Range("G63").Value = 'Sanicula europaea'
$ws.Range("H63").Value = 'L.'
$ws.Range("I63").Value = ''
$ws.Range("J63").Value = ''
$ws.Range("L63").Value = ''
$ws.Range("P63").Value = 'Fiskarsundet, Srm'
$ws.Range("Q63").Value = 689111.5690902721
$ws.Range("R63").Value = 6570305.953062683
$ws.Range("S63").Value = 23
$ws.Range("Z63").Value = '09:34'
$ws.Range("AB63").Value = '09:34'
$ws.Range("A65").Value = 111683853
$ws.Range("B65").Value = 90687
$ws.Range("E65").Value = 5964
$ws.Range("F65").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G65").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H65").Value = '(L.:Fr.) P.Karst.'
$ws.Range("L65").ClearContents()
$ws.Range("A66").Value = 111683850
$ws.Range("B66").Value = 90332
$ws.Range("E66").Value = 4769
$ws.Range("F66").Value = 'Svavelriska'
$ws.Range("G66").Value = 'Lactarius scrobiculatus'
$ws.Range("H66").Value = '(Scop.:Fr.) Fr.'
$ws.Range("I66").Value = '''3'
$ws.Range("J66").Value = 'fruktkroppar'
$ws.Range("P66").Value = 'Bergaholm, Tyresö kn, Srm'
$ws.Range("Q66").Value = 689075.4602011892
$ws.Range("R66").Value = 6570319.534944151
$ws.Range("S66").Value = 20
$ws.Range("Z66").Value = '09:25'
$ws.Range("AB66").Value = '09:25'
